$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.642.48'
$ws.Range('E2').Value = '  -0.66%  '

$ws.Range('D3').Value = '2.453.40'
$ws.Range('E3').Value = '  -0.68%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.43%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -1.87%  '

$ws.Range('E9').Value = '  -1.33%  '

$ws.Range('E10').Value = '  -0.23%  '

$ws.Range('E11').Value = '  -2.20%  '

$ws.Range('E12').Value = '  -1.73%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.74%  '

$ws.Range('E14').Value = '  -3.59%  '

$ws.Range('D15').Value = '2.900.21'
$ws.Range('E15').Value = '  -0.63%  '

$ws.Range('D16').Value = '62.535.78'
$ws.Range('E16').Value = '  -0.86%  '

$ws.Range('D17').Value = '2.455.65'
$ws.Range('E17').Value = '  -0.62%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.73%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.05%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '321.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.49%  '

$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.64%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.29%  '

$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.41%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '644.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.57%  '

$ws.Range('D27').Value = '2.580.20'
$ws.Range('E27').Value = '  -0.50%  '

$ws.Range('D28').Value = '0.0₃0949'
$ws.Range('E28').Value = '  -4.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.16%  '

$ws.Range('E30').Value = '  -3.50%  '

$ws.Range('E31').Value = '  -2.78%  '

$ws.Range('E32').Value = '  -3.34%  '

$ws.Range('E33').Value = '  -0.53%  '

$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '150.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.11%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.362'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.47%  '

$ws.Range('E40').Value = '  -3.03%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.11%  '

$ws.Range('E43').Value = '  +1.99%  '

$ws.Range('E44').Value = '  +0.64%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '152.47'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('E46').Value = '  +1.75%  '

$ws.Range('E47').Value = '  -2.01%  '

$ws.Range('E48').Value = '  -0.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.42%  '

$ws.Range('E50').Value = '  -1.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0900'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.98%  '
